# Hybrid parameter passing feature:
# - "Search on UI page" step now sends the {{empId}} parameter (exported
#   earlier from the customer API response) instead of the hard-coded
#   literal "seattle".
# - A new assertion override is added that checks the page title against
#   the same {{empId}} parameter.
# - The button-click step now targets the resolved locator
#   "Google.googlepage.searchbutton" directly instead of the unused
#   {{searchBar}} placeholder.
# - The googletest sheet becomes the active sheet/tab in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("t-googletest")

# Row 8: sendKeys step - use the {{empId}} parameter instead of "seattle",
# and add an Overrides assertion that checks the page title.
$ws.Range("D8").Value = "{{empId}}"
$ws.Range("E8").Value = "assertEquals::WebDriver::getTitle::{{empId}}"

# Row 9: click step - target resolved locator instead of {{searchBar}}.
$ws.Range("C9").Value = "Google.googlepage.searchbutton"

# Column E got a bit narrower once laid out with the new content.
$ws.Columns.Item(5).ColumnWidth = 48.6

# Make the googletest sheet the active tab (activeTab goes from 3 to 1,
# tabSelected moves from d-dev to t-googletest).
$ws.Activate()
